# Auto-generated edit script: applies updated Profits values per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 7800.222
$ws.Range("I40").Value = 7467.222
$ws.Range("J40").Value = 8133.222
$ws.Range("K40").Value = 7467.222
$ws.Range("L40").Value = 8133.222
$ws.Range("M40").Value = -7292.222
$ws.Range("N40").Value = -8483.222

# Row 70
$ws.Range("H70").Value = 60736.883
$ws.Range("I70").Value = 912.5
$ws.Range("J70").Value = 93368.37
$ws.Range("K70").Value = 2737.5
$ws.Range("L70").Value = 280105.11
$ws.Range("M70").Value = -2467.5
$ws.Range("N70").Value = -280645.11

# Row 73
$ws.Range("H73").Value = 60736.883
$ws.Range("I73").Value = 912.5
$ws.Range("J73").Value = 93368.37
$ws.Range("K73").Value = 2737.5
$ws.Range("L73").Value = 280105.11
$ws.Range("M73").Value = -1801.5
$ws.Range("N73").Value = -281977.11

# Row 137
$ws.Range("H137").Value = 516196.22
$ws.Range("I137").Value = 418844.47
$ws.Range("J137").Value = 671959.0600000001
$ws.Range("K137").Value = 1256533.41
$ws.Range("L137").Value = 2015877.18
$ws.Range("M137").Value = -1253983.41
$ws.Range("N137").Value = -2020977.18

# Row 138
$ws.Range("H138").Value = 3791.492
$ws.Range("I138").Value = 2106.5625
$ws.Range("J138").Value = 9183.267
$ws.Range("K138").Value = 6319.6875
$ws.Range("L138").Value = 27549.801
$ws.Range("M138").Value = -1179.6875
$ws.Range("N138").Value = -37829.801


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 232
$ws.Range("I5").Value = 232
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 232
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -120

# Row 32
$ws.Range("H32").Value = 3484.125
$ws.Range("I32").Value = 2983.6128
$ws.Range("J32").Value = 19000
$ws.Range("K32").Value = 2983.6128
$ws.Range("L32").Value = 19000
$ws.Range("M32").Value = -2696.6128
$ws.Range("N32").Value = -19574

# Row 45
$ws.Range("H45").Value = 6235.154
$ws.Range("I45").Value = 5256.5
$ws.Range("J45").Value = 7801
$ws.Range("K45").Value = 5256.5
$ws.Range("L45").Value = 7801
$ws.Range("M45").Value = -4879.5
$ws.Range("N45").Value = -8555

# Row 74
$ws.Range("H74").Value = 1590.0333
$ws.Range("I74").Value = 1590.0769
$ws.Range("J74").Value = 1589.75
$ws.Range("K74").Value = 1590.0769
$ws.Range("L74").Value = 1589.75
$ws.Range("M74").Value = -716.0769

# Row 77
$ws.Range("H77").Value = 1590.0333
$ws.Range("I77").Value = 1590.0769
$ws.Range("J77").Value = 1589.75
$ws.Range("K77").Value = 7950.3845
$ws.Range("L77").Value = 7948.75
$ws.Range("M77").Value = -3582.3845

# Row 122
$ws.Range("H122").Value = 3946.7307
$ws.Range("I122").Value = 2054.6365
$ws.Range("J122").Value = 5334.2666
$ws.Range("K122").Value = 6163.9095
$ws.Range("L122").Value = 16002.7998
$ws.Range("M122").Value = -3713.9095


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 232
$ws.Range("I4").Value = 232
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 232
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -117

# Row 70
$ws.Range("H70").Value = 186000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 186000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 186000
$ws.Range("N70").Value = -186586

# Row 73
$ws.Range("H73").Value = 186000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 186000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 186000
$ws.Range("N73").Value = -188028

# Row 107
$ws.Range("H107").Value = 591046.3
$ws.Range("I107").Value = 2841
$ws.Range("J107").Value = 3336004.2
$ws.Range("K107").Value = 2841
$ws.Range("L107").Value = 3336004.2
$ws.Range("M107").Value = -921


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 265936.9
$ws.Range("I31").Value = 304803.62
$ws.Range("J31").Value = 9416.4
$ws.Range("K31").Value = 304803.62
$ws.Range("L31").Value = 9416.4
$ws.Range("M31").Value = -304508.62
$ws.Range("N31").Value = -10006.4

# Row 34
$ws.Range("H34").Value = 265936.9
$ws.Range("I34").Value = 304803.62
$ws.Range("J34").Value = 9416.4
$ws.Range("K34").Value = 304803.62
$ws.Range("L34").Value = 9416.4
$ws.Range("M34").Value = -304601.62
$ws.Range("N34").Value = -9820.4

# Row 109
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 132
$ws.Range("H132").Value = 2660.38
$ws.Range("I132").Value = 1914
$ws.Range("J132").Value = 4109.2354
$ws.Range("K132").Value = 5742
$ws.Range("L132").Value = 12327.7062
$ws.Range("M132").Value = -3212

# Row 134
$ws.Range("H134").Value = 284604.94
$ws.Range("I134").Value = 180546.53
$ws.Range("J134").Value = 673089.6
$ws.Range("K134").Value = 541639.59
$ws.Range("L134").Value = 2019268.8
$ws.Range("M134").Value = -539104.59


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 3515
$ws.Range("I2").Value = 4322.231
$ws.Range("J2").Value = 17
$ws.Range("K2").Value = 25933.386
$ws.Range("L2").Value = 102
$ws.Range("M2").Value = -25820.386
$ws.Range("N2").Value = -328

# Row 12
$ws.Range("H12").Value = 1088.7
$ws.Range("I12").Value = 19.5
$ws.Range("J12").Value = 1356
$ws.Range("K12").Value = 58.5
$ws.Range("L12").Value = 4068
$ws.Range("M12").Value = 114.5
$ws.Range("N12").Value = -4414

# Row 32
$ws.Range("H32").Value = 4803495.5
$ws.Range("I32").Value = 7500072.5
$ws.Range("J32").Value = 3005777.5
$ws.Range("K32").Value = 22500217.5
$ws.Range("L32").Value = 9017332.5
$ws.Range("M32").Value = -22499934.5

# Row 92
$ws.Range("H92").Value = 556204.4
$ws.Range("I92").Value = 1428984.6
$ws.Range("J92").Value = 798.8182
$ws.Range("K92").Value = 4286953.800000001
$ws.Range("L92").Value = 2396.4546
$ws.Range("M92").Value = -4285705.800000001


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2446.516
$ws.Range("I102").Value = 1572.579
$ws.Range("J102").Value = 3830.25
$ws.Range("K102").Value = 1572.579
$ws.Range("L102").Value = 3830.25
$ws.Range("M102").Value = 49.42100000000005
$ws.Range("N102").Value = -7074.25

# Row 132
$ws.Range("H132").Value = 480022.38
$ws.Range("I132").Value = 608305.3
$ws.Range("J132").Value = 95173.55
$ws.Range("K132").Value = 1824915.9
$ws.Range("L132").Value = 285520.65
$ws.Range("M132").Value = -1822385.9


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 4983
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 4983
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 4983
$ws.Range("N2").Value = -5207

# Row 68
$ws.Range("H68").Value = 4737.7
$ws.Range("I68").Value = 3214.6667
$ws.Range("J68").Value = 5390.4287
$ws.Range("K68").Value = 3214.6667
$ws.Range("L68").Value = 5390.4287
$ws.Range("M68").Value = -2465.6667
$ws.Range("N68").Value = -6888.4287

# Row 71
$ws.Range("H71").Value = 4737.7
$ws.Range("I71").Value = 3214.6667
$ws.Range("J71").Value = 5390.4287
$ws.Range("K71").Value = 16073.3335
$ws.Range("L71").Value = 26952.1435
$ws.Range("M71").Value = -12329.3335
$ws.Range("N71").Value = -34440.14350000001

# Row 74
$ws.Range("H74").Value = 36000.4
$ws.Range("I74").Value = 50197
$ws.Range("J74").Value = 34423
$ws.Range("K74").Value = 50197
$ws.Range("L74").Value = 34423
$ws.Range("M74").Value = -49199
$ws.Range("N74").Value = -36419

# Row 77
$ws.Range("H77").Value = 36000.4
$ws.Range("I77").Value = 50197
$ws.Range("J77").Value = 34423
$ws.Range("K77").Value = 150591
$ws.Range("L77").Value = 103269
$ws.Range("M77").Value = -145599
$ws.Range("N77").Value = -113253

# Row 108
$ws.Range("H108").Value = 40010
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 40010
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 40010
$ws.Range("N108").Value = -47690

# Row 122
$ws.Range("H122").Value = 317097.88
$ws.Range("I122").Value = 458520.5
$ws.Range("J122").Value = 5968.1
$ws.Range("K122").Value = 1375561.5
$ws.Range("L122").Value = 17904.3
$ws.Range("M122").Value = -1373111.5

# Row 132
$ws.Range("H132").Value = 4318.148
$ws.Range("I132").Value = 2997.3142
$ws.Range("J132").Value = 6751.263
$ws.Range("K132").Value = 8991.942599999998
$ws.Range("L132").Value = 20253.789
$ws.Range("M132").Value = -6461.942599999998

# Row 136
$ws.Range("H136").Value = 235248.89
$ws.Range("I136").Value = 379367.16
$ws.Range("J136").Value = 3786.2122
$ws.Range("K136").Value = 1138101.48
$ws.Range("L136").Value = 11358.6366
$ws.Range("M136").Value = -1135551.48


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 500.3684
$ws.Range("I107").Value = 518
$ws.Range("J107").Value = 462.16666
$ws.Range("K107").Value = 1554
$ws.Range("L107").Value = 1386.49998
$ws.Range("M107").Value = 366

# Row 122
$ws.Range("H122").Value = 28574094
$ws.Range("I122").Value = 41668110
$ws.Range("J122").Value = 5329.4546
$ws.Range("K122").Value = 125004330
$ws.Range("L122").Value = 15988.3638
$ws.Range("M122").Value = -125001880

# Row 136
$ws.Range("H136").Value = 227147.92
$ws.Range("I136").Value = 261006.7
$ws.Range("J136").Value = 136857.86
$ws.Range("K136").Value = 783020.1000000001
$ws.Range("L136").Value = 410573.58
$ws.Range("M136").Value = -780470.1000000001

